$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 12
$ws.Range("B3").Value = 125
$ws.Range("C3").Select()
